# Update the "想去人数" (want-to-go count) column F values on three sheets
# to match the newly scraped numbers, as described in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 1154
    3  = 1964
    4  = 619
    5  = 1273
    6  = 71
    7  = 50
    9  = 344
    10 = 128
    11 = 104
    12 = 855
    13 = 262
    14 = 136
    17 = 348
    18 = 254
    19 = 708
    20 = 83
    21 = 673
    22 = 203
    24 = 916
    25 = 374
    26 = 200
    27 = 59
    31 = 430
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    4  = 335
    5  = 21
    6  = 32
    11 = 131
    12 = 27
}
foreach ($row in $updates2.Keys) {
    $ws2.Range("F$row").Value = $updates2[$row]
}

# Sheet "全部类型" (All types, combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    3  = 1154
    4  = 1964
    5  = 619
    6  = 1273
    7  = 71
    9  = 50
    11 = 344
    12 = 128
    13 = 104
    14 = 855
    15 = 262
    16 = 136
    19 = 335
    21 = 21
    22 = 348
    23 = 32
    25 = 254
    26 = 708
    27 = 83
    28 = 673
    29 = 203
    31 = 916
    32 = 374
    35 = 200
    36 = 59
    39 = 131
    42 = 27
    43 = 430
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
